# ==========================================================================
# PlayerPerformance_3465.xlsx edit:
#   1. Insert a new "Player Info" sheet at the front with player bio data.
#   2. Rename the MATCH_CARD_LINK column to MATCH_CODE on "ODI Batting" and
#      "ODI Bowling", replacing each URL value with just the numeric
#      MatchCode query-param.
#   3. Append a new "ODI Batting Extra" sheet at the end with additional
#      per-match batting stats.
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------------
# Helper: apply the bold / thin-border / center-top-aligned header look
# used by row 1 on the existing sheets to a freshly written header range.
# --------------------------------------------------------------------------
function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

# --------------------------------------------------------------------------
# 1. "Player Info" sheet (becomes the first sheet / sheetId 1)
# --------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Add()
$wsInfo.Name = "Player Info"

$wsInfo.Range("A1").Value = "ID"
$wsInfo.Range("B1").Value = "NAME"
$wsInfo.Range("C1").Value = "BATTING_HAND"
$wsInfo.Range("D1").Value = "BOWL_STYLE"
Set-HeaderStyle $wsInfo.Range("A1:D1")

$wsInfo.Range("A2").Value = "'3465"
$wsInfo.Range("B2").Value = "Kieron Adrian Pollard"
$wsInfo.Range("C2").Value = "Right Handed"
$wsInfo.Range("D2").Value = "Right Arm Medium"

# --------------------------------------------------------------------------
# 2a. "ODI Batting" sheet: MATCH_CARD_LINK (col D) -> MATCH_CODE
# --------------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"

$lastRowBatting = $wsBatting.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowBatting; $r++) {
    $cell = $wsBatting.Cells.Item($r, 4)
    $link = $cell.Value()
    if ($link -ne $null -and $link -ne "") {
        $code = $link -replace ".*MatchCode=", ""
        $cell.Value = "'" + $code
    }
}

# --------------------------------------------------------------------------
# 2b. "ODI Bowling" sheet: MATCH_CARD_LINK (col B) -> MATCH_CODE
# --------------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsBowling.Cells.Item(1, 2).Value = "MATCH_CODE"

$lastRowBowling = $wsBowling.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowBowling; $r++) {
    $cell = $wsBowling.Cells.Item($r, 2)
    $link = $cell.Value()
    if ($link -ne $null -and $link -ne "") {
        $code = $link -replace ".*MatchCode=", ""
        $cell.Value = "'" + $code
    }
}

# --------------------------------------------------------------------------
# 3. "ODI Batting Extra" sheet (new, appended at the end)
# --------------------------------------------------------------------------
$wsExtra = $wb.Worksheets.Add()
$wsExtra.Name = "ODI Batting Extra"
$wsExtra.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"
Set-HeaderStyle $wsExtra.Range("A1:F1")

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4379", "6", "1", "2", "12.65%", "NO"),
    @("4385", "",  "",  "",  "",       "NO"),
    @("4387", "6", "0", "0", "",       "NO"),
    @("4388", "6", "3", "7", "23.49%", "NO"),
    @("4391", "7", "0", "0", "",       "NO"),
    @("4394", "6", "2", "4", "16.53%", "NO"),
    @("4397", "8", "",  "",  "",       "NO"),
    @("4413", "6", "0", "1", "3.11%",  "NO"),
    @("4414", "",  "",  "",  "",       "NO"),
    @("4417", "4", "6", "1", "16.28%", "NO"),
    @("4449", "6", "",  "",  "",       "NO"),
    @("4450", "",  "",  "",  "",       "NO"),
    @("4451", "6", "4", "1", "19.20%", "NO"),
    @("4483", "6", "5", "3", "45.53%", "NO"),
    @("4484", "6", "0", "0", "1.05%",  "NO"),
    @("4486", "",  "",  "",  "",       "NO"),
    @("4519", "",  "",  "",  "",       "NO"),
    @("4520", "6", "0", "0", "0.44%",  "NO"),
    @("4522", "",  "",  "",  "",       "NO"),
    @("4533", "",  "",  "",  "",       "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $wsExtra.Cells.Item($r, 1).Value = "'" + $row[0]

    if ($row[1] -ne "") {
        $wsExtra.Cells.Item($r, 2).Value = [double]$row[1]
    } else {
        $wsExtra.Cells.Item($r, 2).Value = ""
    }

    if ($row[2] -ne "") {
        $wsExtra.Cells.Item($r, 3).Value = "'" + $row[2]
    } else {
        $wsExtra.Cells.Item($r, 3).Value = ""
    }

    if ($row[3] -ne "") {
        $wsExtra.Cells.Item($r, 4).Value = "'" + $row[3]
    } else {
        $wsExtra.Cells.Item($r, 4).Value = ""
    }

    if ($row[4] -ne "") {
        $wsExtra.Cells.Item($r, 5).Value = "'" + $row[4]
    } else {
        $wsExtra.Cells.Item($r, 5).Value = ""
    }

    $wsExtra.Cells.Item($r, 6).Value = "'" + $row[5]

    $r = $r + 1
}
